# Applies the crypto-price/volume refresh from the Tue Jan 3 17:35:54 UTC 2023
# "Updated symbol list" GitHub Actions commit.
#
# The source cells are plain text (t="inlineStr") even though many values
# look numeric (prices, percentages). Setting .Value directly on a
# numeric-looking string would make Excel COM silently re-type the cell as
# a Number (and introduce float rounding noise, e.g. 245.26 -> 245.25999...).
# To keep the cells as text we briefly force NumberFormat "@" (Text) before
# assigning the value, then restore the "Normal" style so no stray
# number-format style lingers on the cell afterward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "245.26"
Set-TextValue "E2" "-0.49%"
Set-TextValue "D3" "29.19"
Set-TextValue "E3" "-1.09%"
Set-TextValue "D4" "5.264"
Set-TextValue "D5" "0.05706"
Set-TextValue "E5" "-0.01%"
Set-TextValue "D6" "6.615"
Set-TextValue "E6" "0.22%"
Set-TextValue "D7" "3.194"
Set-TextValue "E7" "3.96%"
Set-TextValue "D8" "0.8513"
Set-TextValue "E8" "-0.75%"
Set-TextValue "D9" "0.8571"
Set-TextValue "E9" "-1.57%"
Set-TextValue "D10" "0.1369"
Set-TextValue "E10" "0.23%"
Set-TextValue "D11" "0.07048"
Set-TextValue "E11" "-0.64%"
Set-TextValue "D12" "0.03191"
Set-TextValue "E12" "9.37%"
Set-TextValue "D13" "0.09286"
Set-TextValue "E13" "-1.13%"
Set-TextValue "D14" "0.001526"
Set-TextValue "E14" "0.03%"
Set-TextValue "D15" "0.0005983"
Set-TextValue "E15" "-0.11%"
Set-TextValue "D16" "0.005928"
Set-TextValue "E16" "-0.76%"
Set-TextValue "E17" "0.25%"
Set-TextValue "E18" "-0.49%"
Set-TextValue "D19" "0.3161"
Set-TextValue "E19" "-0.42%"
Set-TextValue "D20" "0.03315"
Set-TextValue "E20" "-0.25%"
Set-TextValue "E21" "-1.91%"
Set-TextValue "D22" "3.508"
Set-TextValue "E22" "0.93%"
Set-TextValue "D23" "0.04101"
Set-TextValue "E23" "-1.94%"
Set-TextValue "D24" "0.1380"
Set-TextValue "E24" "-0.01%"
Set-TextValue "D25" "0.001225"
Set-TextValue "E25" "0.29%"
Set-TextValue "D26" "0.004143"
Set-TextValue "E26" "-17.64%"
Set-TextValue "E27" "-0.80%"
Set-TextValue "E28" "-25.22%"
Set-TextValue "D40" "0.03754"
Set-TextValue "E40" "0.27%"
Set-TextValue "E41" "-0.90%"
Set-TextValue "D42" "0.003712"
Set-TextValue "E42" "6.65%"
Set-TextValue "E43" "-3.59%"
Set-TextValue "D44" "0.009347"
Set-TextValue "E44" "-6.06%"
Set-TextValue "D45" "0.00005273"
Set-TextValue "E45" "1.13%"
Set-TextValue "E46" "0.06%"
Set-TextValue "E47" "29.39%"
Set-TextValue "D48" "0.002442"
Set-TextValue "E48" "-4.74%"
Set-TextValue "E49" "0.06%"
Set-TextValue "E50" "0.06%"
